# Add a new scene quest row to the Drop table, and rename the
# "suijilv#zhuangbei" Ename codes to "zhuangbei#".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- 1. Fill in the new row (row 41) ------------------------------------
# Match the style/formatting of the row above it (row 40) first, the same
# way Excel copies formatting down when a new table row is authored.
$ws.Range("B40:E40").Copy()
$ws.Range("B41:E41").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I40").Copy()
$ws.Range("I41").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = 0

# Column order as originally authored: B (Name), D (Items), E (ItemRate),
# C (Ename), A (Id), I (Count).
$ws.Cells.Item(41, 2).Value = "sq隐秘石门"
$ws.Cells.Item(41, 4).Value = "sucaidai;kapaibugeibao(wu);xiaoxingmofayaoji;jingyanzhishu;suijihuanshouka"
$ws.Cells.Item(41, 5).Value = "15;20;20;20;25"
$ws.Cells.Item(41, 3).Value = "dlshimen"
$ws.Cells.Item(41, 1).Value = 23000501
$ws.Cells.Item(41, 9).Value = 2

# --- 2. Rename the Ename codes for the 5 random-equipment rows ----------
# (Order matters for shared-string ordering: C10, C11, C12, C13, then C9.)
$ws.Cells.Item(10, 3).Value = "dlzhuangbei2"
$ws.Cells.Item(11, 3).Value = "dlzhuangbei3"
$ws.Cells.Item(12, 3).Value = "dlzhuangbei4"
$ws.Cells.Item(13, 3).Value = "dlzhuangbei5"
$ws.Cells.Item(9, 3).Value = "dlzhuangbei1"

# --- 3. Expand the table / autofilter range to include the new row ------
$tbl.Resize($ws.Range("A3:I41"))

# --- 4. Update the active selection, like the author's last click -------
$null = $ws.Range("D15").Select()
